# Insert two new weekly price rows for "Pimiento" (Zafiro rojo / Zafiro verde)
# at the top of the date-ordered block that starts at row 117, pushing the
# existing rows 117-131 down to 119-133.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 117:131 down by two rows (new blank rows appear at 117:118)
$ws.Rows("117:118").Insert()

# --- New row 117: Zafiro rojo ---
$ws.Range("A117").Value = 7
$ws.Range("B117").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C117").Value = "Ñuble"
$ws.Range("D117").Value = 44449
$ws.Range("E117").Value = 16
$ws.Range("F117").Value = 100112002
$ws.Range("G117").Value = "Pimiento"
$ws.Range("H117").Value = "Zafiro rojo"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 160
$ws.Range("K117").Value = 44000
$ws.Range("L117").Value = 45000
$ws.Range("M117").Value = 44500
$ws.Range("N117").Value = "$/caja 15 kilos"
$ws.Range("O117").Value = "Región de Arica y Parinacota"
$ws.Range("P117").Value = 2967
$ws.Range("Q117").Value = 15
$ws.Range("R117").Value = "Hortaliza"

# --- New row 118: Zafiro verde ---
$ws.Range("A118").Value = 7
$ws.Range("B118").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C118").Value = "Ñuble"
$ws.Range("D118").Value = 44449
$ws.Range("E118").Value = 16
$ws.Range("F118").Value = 100112002
$ws.Range("G118").Value = "Pimiento"
$ws.Range("H118").Value = "Zafiro verde"
$ws.Range("I118").Value = "Primera"
$ws.Range("J118").Value = 160
$ws.Range("K118").Value = 41000
$ws.Range("L118").Value = 42000
$ws.Range("M118").Value = 41500
$ws.Range("N118").Value = "$/caja 15 kilos"
$ws.Range("O118").Value = "Región de Arica y Parinacota"
$ws.Range("P118").Value = 2767
$ws.Range("Q118").Value = 15
$ws.Range("R118").Value = "Hortaliza"
